# Adds new "support scripts" rows (auto-start + RTC) to the Commands sheet,
# plus a new column C (and trailing column D filler) of single-space cells,
# mirroring the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- New content -----------------------------------------------------------
# Row 3: B3 gets a hyperlink to the auto-start guide (new shared string #14)
$ws.Range("B3").Value = "http://www.opentechguides.com/how-to/article/raspberry-pi/5/raspberry-pi-auto-start.html"

# Row 9 (new row): RTC entry, with both an adafruit link (B9) and a repeat
# of the auto-start link (C9)
$ws.Range("A9").Value = "RTC"
$ws.Range("B9").Value = "https://learn.adafruit.com/adding-a-real-time-clock-to-raspberry-pi?view=all"

# Fill column C (rows 2-8) and column D (rows 2-14) with a single space,
# matching the new shared string entry added at the end (" ").
$ws.Range("C2:C8").Value = " "
$ws.Range("D2:D14").Value = " "

# C9 repeats the auto-start URL text (reuses shared string #14)
$ws.Range("C9").Value = "http://www.opentechguides.com/how-to/article/raspberry-pi/5/raspberry-pi-auto-start.html"

# --- Column width for the new column C --------------------------------------
$ws.Columns.Item(3).ColumnWidth = 27.86

# --- Hyperlinks --------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B3"), "http://www.opentechguides.com/how-to/article/raspberry-pi/5/raspberry-pi-auto-start.html")
$ws.Hyperlinks.Add($ws.Range("C9"), "http://www.opentechguides.com/how-to/article/raspberry-pi/5/raspberry-pi-auto-start.html")
$ws.Hyperlinks.Add($ws.Range("B9"), "https://learn.adafruit.com/adding-a-real-time-clock-to-raspberry-pi?view=all")

# Give the new hyperlinked cells the same "Hyperlink" style used elsewhere
$ws.Range("B3").Style = $ws.Range("B2").Style
$ws.Range("B9").Style = $ws.Range("B2").Style
$ws.Range("C9").Style = $ws.Range("B2").Style

# --- Selection (matches the saved view state in the authored workbook) -----
$ws.Range("C15").Select() | Out-Null
